# Apply "add regular/private online training format" edit requested by restu.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Try to (re)apply codenames; harmless if the runtime treats these as read-only.
$wb.CodeName = "ThisWorkbook"
$ws.CodeName = "Sheet1"

# New block describing the Online Training type (rows 30-35). Populated in
# this order so new shared-string entries line up with the source workbook.
$ws.Range("A30").Value = "Tipe Online Training:"
$ws.Range("B30").Value = "ROT"

$ws.Range("A31").Value = "(Pilih Salah Satu, ROT, POT atau 1H)"
# B31 stays empty, but picks up the centered placeholder style used elsewhere
# in column A/B for blank cells awaiting input (same style as A11).
$ws.Range("A11").Copy() | Out-Null
$ws.Range("B31").PasteSpecial(-4122) | Out-Null

$ws.Range("A33").Value = "ROT => Regular Online Training"
$ws.Range("A34").Value = "POT => Private Online Training"
$ws.Range("A35").Value = "1H => 1Hour Online Training"

# A32 is a bold section divider, matching the style already used for the
# "--------------------" divider at A10. The leading apostrophe forces the
# "===" text to be stored literally instead of being parsed as a formula.
$ws.Range("A10").Copy() | Out-Null
$ws.Range("A32").PasteSpecial(-4122) | Out-Null
$ws.Range("A32").Value = "'==="

# New "Meeting ID" sample value next to the existing label.
$ws.Range("B24").Value = "11222 3344"

# The training name shown near the top changed from an "Inhouse" test entry
# to a "Private Online Training" test entry.
$ws.Range("B2").Value = "Private Online Training Testing Brainmatics Training"

# Update the active selection to C7, as recorded in the saved workbook.
$ws.Range("C7").Select() | Out-Null

$excel.CutCopyMode = $false
